# Refresh cryptos list snapshot (prices + 1h change %) - GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.310.34"
$ws.Range("E2").Value = "  +2.92%  "

$ws.Range("D3").Value = "'2.413.26"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'572.22"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").Value = "'144.28"
$ws.Range("E6").Value = "  +4.58%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'2.438.79"
$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("E10").Value = "  +4.60%  "

$ws.Range("E11").Value = "  +0.64%  "

$ws.Range("D12").Value = "'5.25"

$ws.Range("E13").Value = "  +3.89%  "

$ws.Range("D14").Value = "'26.85"
$ws.Range("E14").Value = "  +4.61%  "

$ws.Range("D15").Value = "'0.0000180"
$ws.Range("E15").Value = "  +8.60%  "

$ws.Range("D16").Value = "'2.870.13"
$ws.Range("E16").Value = "  +2.01%  "

$ws.Range("D17").Value = "'62.119.31"
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("D18").Value = "'2.423.33"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "'7.90"
$ws.Range("E19").Value = "  -4.45%  "

$ws.Range("D20").Value = "'10.88"
$ws.Range("E20").Value = "  +2.64%  "

$ws.Range("D21").Value = "'325.89"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("E22").Value = "  +2.38%  "

$ws.Range("E23").Value = "  +13.19%  "

$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "'65.46"
$ws.Range("E25").Value = "  +1.70%  "

$ws.Range("D26").Value = "'613.77"
$ws.Range("E26").Value = "  +10.42%  "

$ws.Range("D27").Value = "'8.38"
$ws.Range("E27").Value = "  +4.38%  "

$ws.Range("D28").Value = "'0.0₃0984"
$ws.Range("E28").Value = "  +8.05%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "'2.557.63"
$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'0.994"
$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("E31").Value = "  +2.32%  "

$ws.Range("E32").Value = "  +8.41%  "

$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.136"
$ws.Range("E33").Value = "  +4.85%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  +1.86%  "

$ws.Range("E35").Value = "  +4.82%  "

$ws.Range("D36").Value = "'0.995"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("E37").Value = "  +5.39%  "

$ws.Range("D38").Value = "'152.74"
$ws.Range("E38").Value = "  -0.33%  "

$ws.Range("D39").Value = "'0.372"
$ws.Range("E39").Value = "  +1.31%  "

$ws.Range("D40").Value = "'5.39"
$ws.Range("E40").Value = "  +6.41%  "

$ws.Range("D41").Value = "'18.57"
$ws.Range("E41").Value = "  +1.75%  "

$ws.Range("D42").Value = "'2.69"
$ws.Range("E42").Value = "  +17.01%  "

$ws.Range("E43").Value = "  +6.63%  "

$ws.Range("D44").Value = "'42.33"
$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").Value = "'0.0₆0280"
$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("D47").Value = "'143.40"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").Value = "'3.58"
$ws.Range("E48").Value = "  +2.59%  "

$ws.Range("D49").Value = "'20.28"
$ws.Range("E49").Value = "  +7.04%  "

$ws.Range("D50").Value = "'0.601"
$ws.Range("E50").Value = "  +2.14%  "

$ws.Range("D51").Value = "'0.0513"
$ws.Range("E51").Value = "  +3.27%  "
